$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old header row (row 1). Using ClearContents (rather than
# Delete) keeps every other row at its original row number - row 2 stays
# row 2, row 9 stays row 9, etc. - matching the target sheet where the
# used range simply starts at row 2 afterwards.
$ws.Rows.Item(1).ClearContents()

# Append a new faculty record as row 10.
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = "newfaculty@university.edu"
$ws.Range("D10").Value = "newpass"
$ws.Range("E10").Value = "Faculty"
$ws.Range("F10").Value = $false

# G10/H10 mirror the blank-string placeholder used by the other generated
# rows (e.g. G6/H6): a single quote forces Excel to store an empty text
# value in the cell instead of leaving it truly blank.
$ws.Range("G10").Value = "'"
$ws.Range("H10").Value = "'"
$ws.Range("G10").Style = "Normal"
$ws.Range("H10").Style = "Normal"

$ws.Range("I10").Value = 0
